$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4
    $ws.Range("F4").Value = 1441
    $ws.Range("F6").Value = 23
    $ws.Range("F9").Value = 228
}
